$wb = $excel.ActiveWorkbook

# --- sheet "participants": rename two header labels ---
$ws1 = $wb.Worksheets.Item("participants")
$ws1.Range("H1").Value = "gestational age"
$ws1.Range("J1").Value = "birth weight"

# --- sheet "sessions": insert a new column "task_positions" right after "tasks" (col O) ---
$ws2 = $wb.Worksheets.Item("sessions")
$ws2.Columns.Item(16).Insert()
$ws2.Range("P1").Value = "task_positions"
